$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) holds text, but some updated values (e.g. "1.002")
# would otherwise be auto-detected as numbers by Excel. Force the whole
# Price column to stay formatted as Text before writing the new values,
# matching the original inline-string (non-numeric) cell content.
$ws.Range("D2:D51").NumberFormat = "@"

# --- Updated crypto list values ---
$ws.Range("D2").Value = "23.217.67"
$ws.Range("E2").Value = "  +0.67%  "
$ws.Range("D3").Value = "1.604.37"
$ws.Range("E3").Value = "  +0.35%  "
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  +0.20%  "
$ws.Range("D5").Value = "1.001"
$ws.Range("E5").Value = "  +0.11%  "
$ws.Range("D6").Value = "304.36"
$ws.Range("E6").Value = "  +0.81%  "
$ws.Range("D7").Value = "0.3764"
$ws.Range("E7").Value = "  -0.51%  "
$ws.Range("D8").Value = "52.53"
$ws.Range("E8").Value = "  +3.64%  "
$ws.Range("D9").Value = "0.3622"
$ws.Range("E9").Value = "  -0.72%  "
$ws.Range("D10").Value = "1.271"
$ws.Range("E10").Value = "  +1.30%  "
$ws.Range("B11").Value = "BinanceUSD"
$ws.Range("C11").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D11").Value = "1.002"
$ws.Range("E11").Value = "  +0.18%  "
$ws.Range("B12").Value = "Dogecoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D12").Value = "0.08134"
$ws.Range("E12").Value = "  -0.07%  "
$ws.Range("D13").Value = "22.89"
$ws.Range("E13").Value = "  +2.20%  "
$ws.Range("D14").Value = "6.593"
$ws.Range("E14").Value = "  +0.02%  "
$ws.Range("B15").Value = "ShibaInu"
$ws.Range("C15").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D15").Value = "0.00001248"
$ws.Range("E15").Value = "  +0.18%  "
$ws.Range("B16").Value = "Chainlink"
$ws.Range("C16").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D16").Value = "7.357"
$ws.Range("E16").Value = "  -0.28%  "
$ws.Range("D17").Value = "1.599.69"
$ws.Range("E17").Value = "  -0.02%  "
$ws.Range("D18").Value = "94.03"
$ws.Range("E18").Value = "  +2.45%  "
$ws.Range("D19").Value = "0.06929"
$ws.Range("E19").Value = "  +1.73%  "
$ws.Range("D20").Value = "18.11"
$ws.Range("E20").Value = "  -0.48%  "
$ws.Range("D21").Value = "6.531"
$ws.Range("E21").Value = "  +0.02%  "
$ws.Range("D22").Value = "1.003"
$ws.Range("E22").Value = "  +0.19%  "
$ws.Range("D23").Value = "12.91"
$ws.Range("E23").Value = "  -0.75%  "
$ws.Range("D24").Value = "23.210.82"
$ws.Range("E24").Value = "  +0.68%  "
$ws.Range("D25").Value = "2.441"
$ws.Range("E25").Value = "  +3.08%  "
$ws.Range("D26").Value = "3.068"
$ws.Range("E26").Value = "  +10.14%  "
$ws.Range("D27").Value = "21.18"
$ws.Range("E27").Value = "  +0.52%  "
$ws.Range("D28").Value = "150.10"
$ws.Range("E28").Value = "  +0.72%  "
$ws.Range("D29").Value = "5.282"
$ws.Range("E29").Value = "  +0.81%  "
$ws.Range("D30").Value = "134.77"
$ws.Range("E30").Value = "  +0.65%  "
$ws.Range("D31").Value = "2.395"
$ws.Range("E31").Value = "  +1.70%  "
$ws.Range("D32").Value = "6.734"
$ws.Range("E32").Value = "  -1.97%  "
$ws.Range("D33").Value = "1.781.01"
$ws.Range("E33").Value = "  +0.29%  "
$ws.Range("D34").Value = "0.9596"
$ws.Range("E34").Value = "  -0.49%  "
$ws.Range("D35").Value = "0.02770"
$ws.Range("E35").Value = "  +2.36%  "
$ws.Range("D36").Value = "0.07444"
$ws.Range("E36").Value = "  -1.96%  "
$ws.Range("D37").Value = "10.36"
$ws.Range("E37").Value = "  +0.35%  "
$ws.Range("D38").Value = "0.2516"
$ws.Range("E38").Value = "  -0.23%  "
$ws.Range("D39").Value = "6.109"
$ws.Range("E39").Value = "  -2.14%  "
$ws.Range("D40").Value = "0.08790"
$ws.Range("E40").Value = "  -0.41%  "
$ws.Range("D41").Value = "1.404"
$ws.Range("E41").Value = "  +2.94%  "
$ws.Range("D42").Value = "0.7086"
$ws.Range("E42").Value = "  +0.54%  "
$ws.Range("D43").Value = "12.40"
$ws.Range("E43").Value = "  +0.37%  "
$ws.Range("D44").Value = "15.83"
$ws.Range("E44").Value = "  +4.12%  "
$ws.Range("D45").Value = "0.6532"
$ws.Range("E45").Value = "  -1.33%  "
$ws.Range("D46").Value = "2.327"
$ws.Range("E46").Value = "  +1.70%  "
$ws.Range("B47").Value = "Frax"
$ws.Range("C47").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D47").Value = "0.9997"
$ws.Range("E47").Value = "  +0.07%  "
$ws.Range("B48").Value = "PancakeSwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D48").Value = "4.009"
$ws.Range("E48").Value = "  +0.40%  "
$ws.Range("B49").Value = "Quant"
$ws.Range("C49").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D49").Value = "133.84"
$ws.Range("E49").Value = "  +1.74%  "
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").Value = "0.07945"
$ws.Range("E50").Value = "  +0.21%  "
$ws.Range("B51").Value = "Flow"
$ws.Range("C51").Value = "https://coinranking.com/coin/QQ0NCmjVq+flow-flow"
$ws.Range("D51").Value = "1.201"
$ws.Range("E51").Value = "  -1.86%  "
